$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (39-46) first, shrinking the sheet
# so the dimension becomes A1:F38.
$ws.Rows("39:46").Delete()

# Column B (ticker) replacements for rows 4-38 (rows 2 and 3 keep their
# original values).
$bValues = @{
  4  = "NSE:ASHAPURMIN"
  5  = "NSE:ATGL"
  6  = "NSE:BALAJITELE"
  7  = "NSE:BHANDARI"
  8  = "NSE:BVCL"
  9  = "NSE:DIVGIITTS"
  10 = "NSE:DJML"
  11 = "NSE:DPSCLTD"
  12 = "NSE:ELGIRUBCO"
  13 = "NSE:EXXARO"
  14 = "NSE:FACT"
  15 = "NSE:FINEORG"
  16 = "NSE:GEECEE"
  17 = "NSE:GENUSPAPER"
  18 = "NSE:GSFC"
  19 = "NSE:HINDCON"
  20 = "NSE:KALYANKJIL"
  21 = "NSE:LICI"
  22 = "NSE:MADRASFERT"
  23 = "NSE:MANGCHEFER"
  24 = "NSE:MANINFRA"
  25 = "NSE:NAGAFERT"
  26 = "NSE:NCLIND"
  27 = "NSE:NFL"
  28 = "NSE:NUVOCO"
  29 = "NSE:ONWARDTEC"
  30 = "NSE:PENIND"
  31 = "NSE:PITTIENG"
  32 = "NSE:PTL"
  33 = "NSE:RADICO"
  34 = "NSE:RAMAPHO"
  35 = "NSE:RELAXO"
  36 = "NSE:RHL"
  37 = "NSE:RPGLIFE"
  38 = "NSE:SAKSOFT"
}
foreach ($row in $bValues.Keys) {
  $ws.Range("B$row").Value = $bValues[$row]
}

# Column C (support zone) replacements for rows 2-4; rows 5-38 are cleared.
$ws.Range("C2").Value = "NSE:BANKA"
$ws.Range("C3").Value = "NSE:MATRIMONY"
$ws.Range("C4").Value = "NSE:RELCHEMQ"
$ws.Range("C5:C38").ClearContents()

# Columns D, E and F no longer hold any data for the remaining rows.
$ws.Range("D2:D38").ClearContents()
$ws.Range("E2:E38").ClearContents()
$ws.Range("F2:F38").ClearContents()
